$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column I entirely (amplitude "init" column), shifting nothing left
# since it is the last used column.
$ws.Range("I1:I14").Delete() | Out-Null

# Update the "amplitude" row (row 9) values to relative ratios instead of
# absolute intensities.
$ws.Range("B9").Value = 0.1
$ws.Range("C9").Value = 0.2
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0.1
$ws.Range("F9").Value = 0.1
$ws.Range("G9").Value = 0.8
$ws.Range("H9").Value = 0.1

# Update the selection to match the recorded workbook state.
$ws.Range("F4").Select() | Out-Null
